# Rename the embedded picture parts in the page header/footers so the
# drawing's display name matches the file it is actually pointing at:
#   - Both Pearson logo pictures (in the two footers) were saved as
#     "image1.png" but should be "image2.png".
#   - The BTEC logo picture (in the header) was saved as "image2.jpg"
#     but should be "image1.jpg".

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footers: Pearson logo, "image1.png" -> "image2.png" ------------------
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

# --- Headers: BTEC logo, "image2.jpg" -> "image1.jpg" ---------------------
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
